$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc = "Several attempts to fix the breach in crest gate number 19 of Tungabhadra dam failed to yield the desired results on Thursday. ADVERTISEMENT."
$flyoverDesc = "The double-decker flyover features two separate carriageways, each with two lanes. It passes by three metro stations: Jayadeva Hospital,..."

# xlVAlignCenter - matches the vertical-center alignment style used by the
# other data rows (cellXfs index 0) so the new/changed cells line up visually
# with the rest of the table.
$xlVAlignCenter = -4108

# Re-assert D12 value and give it the same style as the rest of the table.
$ws.Range("D12").Value = "Outlook India, The Hans India, Deccan Herald"
$ws.Range("D12").VerticalAlignment = $xlVAlignCenter

# Row 13
$ws.Range("A13").Value = "tungabhadra dam gate crash"
$ws.Range("B13").Value = "18 Aug 2024"
$ws.Range("C13").Value = $desc
$ws.Range("D13").Value = "Outlook India, Deccan Herald, The Hans India"
$ws.Range("A13:D13").VerticalAlignment = $xlVAlignCenter

# Row 14
$ws.Range("A14").Value = "double decker flyover in bengaluru"
$ws.Range("B14").Value = "26 Jul 2024"
$ws.Range("C14").Value = $flyoverDesc
$ws.Range("D14").Value = "The Hindu, NDTV, Hindustan Times, Times of India"
$ws.Range("A14:D14").VerticalAlignment = $xlVAlignCenter

# Row 15
$ws.Range("A15").Value = "tungabhadra dam crash"
$ws.Range("B15").Value = "18 Aug 2024"
$ws.Range("C15").Value = $desc
$ws.Range("D15").Value = "Outlook India, The Hans India, Deccan Herald"
$ws.Range("A15:D15").VerticalAlignment = $xlVAlignCenter

# Row 16
$ws.Range("A16").Value = "Tungabhadra Dam crash"
$ws.Range("B16").Value = "18 Aug 2024"
$ws.Range("C16").Value = $desc
$ws.Range("D16").Value = "The Hans India, Deccan Herald, Outlook India"
$ws.Range("A16:D16").VerticalAlignment = $xlVAlignCenter

# Row 17
$ws.Range("A17").Value = "tungabhadra dam gate crash"
$ws.Range("B17").Value = "16 Aug 2024"
$ws.Range("C17").Value = $desc
$ws.Range("D17").Value = "Outlook India, The Hans India, Deccan Herald"
$ws.Range("A17:D17").VerticalAlignment = $xlVAlignCenter

# Row 18
$ws.Range("A18").Value = "tungabhadra dam crash"
$ws.Range("B18").Value = "16 Aug 2024"
$ws.Range("C18").Value = $desc
$ws.Range("D18").Value = "Outlook India, The Hans India, Deccan Herald"
$ws.Range("A18:D18").VerticalAlignment = $xlVAlignCenter
